# Auto-generated edit script: updates cryptos list values per commit
# "Updated cryptos list on Thu Aug 15 15:17:57 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.185.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.06%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.634.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.31%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.06%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'528.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +1.36%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'144.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.72%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.569"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -0.85%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'6.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.55%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +1.17%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.337"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.44%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +0.59%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'3.101.59"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.41%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'59.206.37"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.12%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'20.84"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.29%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'ShibaInu"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0000137"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.53%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'WrappedEther"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'2.670.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.88%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'341.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.11%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +0.78%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'10.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +2.53%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.54%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.40%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'65.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.19%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.417"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +1.16%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +0.03%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.45%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'7.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.05%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.0₃0799"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.45%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "'USDe"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value = "'Aptos"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = "'6.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.90%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'1.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.84%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'18.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.01%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'150.43"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.62%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.06%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -0.01%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.874"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.73%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.863"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.23%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("B38").Value = "'OKB"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'36.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.85%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("B39").Value = "'Stacks"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'1.48"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.42%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'3.64"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +1.32%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.23%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'Stellar"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.0975"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'Mantle"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.601"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -4.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'269.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -2.15%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'19.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.80%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("B46").Value = "'WhiteBITCoin"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'10.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.84%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("B47").Value = "'Hedera"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'0.0537"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.39%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'2.036.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.44%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "'RenderToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'4.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.51%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("B50").Value = "'VeChain"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'0.0229"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.28%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'18.66"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -1.38%  "
$ws.Range("E51").Style = "Normal"
